$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell F1, matching the formatting of the other header cells (B1:E1)
# by copying the format from E1 (bold, centered, bordered) onto F1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Add the time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:52:40.352373"
$ws.Range("F3").Value = "2021-10-05 10:52:40.352384"
$ws.Range("F4").Value = "2021-10-05 10:52:40.352387"
$ws.Range("F5").Value = "2021-10-05 10:52:40.352390"
